## Rent Estimates.xlsx refresh — updated Rentometer/Zillow figures and
## the Rentometer quickview token (-LFNYcE-hBs -> on3TP32HDUw), plus the
## blended rentometer_zillow_user_avg_est numbers that derive from them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Rentometer"
# ---------------------------------------------------------------------
$rentometer = $wb.Worksheets.Item("Rentometer")

$rentometer.Range("B8").Value  = 1242    # average_rent
$rentometer.Range("B9").Value  = 1250    # median_rent
$rentometer.Range("B12").Value = 1200    # percentile_25
$rentometer.Range("B13").Value = 1283    # percentile_75
$rentometer.Range("B14").Value = 62      # std_dev
$rentometer.Range("B15").Value = 13      # sample_size
$rentometer.Range("B18").Value = 1953    # credits_remaining

$newToken = "on3TP32HDUw"

# quickview_url (B17) is a real hyperlink: refresh both the visible text
# and the underlying relationship target, and restore the Hyperlink
# cell style Hyperlinks.Add() just reset.
$newQuickviewUrl = "https://www.rentometer.com/analysis/3-bed/317-newell-st-barberton-oh-44203/$newToken/quickview"
$rentometer.Hyperlinks.Delete()
$rentometer.Hyperlinks.Add($rentometer.Range("B17"), $newQuickviewUrl)
$rentometer.Range("B17").Value = $newQuickviewUrl
$rentometer.Range("B17").Style = "Hyperlink"

$rentometer.Range("B19").Value = $newToken   # token

$rentometer.Range("B20").Value = "[{'rel': 'request pro report', 'href': 'https://www.rentometer.com/api/v1/request_pro_report?api_key=fHSGZM7POi6V5ZPR0w4CXA&token=$newToken'}, {'rel': 'nearby comps', 'href': 'https://www.rentometer.com/api/v1/nearby_comps?api_key=fHSGZM7POi6V5ZPR0w4CXA&token=$newToken'}]"   # links

# ---------------------------------------------------------------------
# Sheet "Zillow"
# ---------------------------------------------------------------------
$zillow = $wb.Worksheets.Item("Zillow")

$zillow.Range("B3").Value  = 1431      # rent_estimate
$zillow.Range("B4").Value  = 870       # min_rent
$zillow.Range("B5").Value  = 1557      # max_rent
$zillow.Range("B8").Value  = 1324.75   # percentile_25
$zillow.Range("B9").Value  = 1413      # percentile_75
$zillow.Range("B10").Value = 1361.5    # median_rent

# ---------------------------------------------------------------------
# Sheet "rentometer_zillow_user_avg_est"
# ---------------------------------------------------------------------
$blend = $wb.Worksheets.Item("rentometer_zillow_user_avg_est")

$blend.Range("B1").Value = 1360.333333333333   # median_rent
$blend.Range("B2").Value = 1357.666666666667   # average_rent
$blend.Range("B3").Value = 1262.375            # percentile_25
$blend.Range("B4").Value = 1348                # percentile_75
